$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.942.45'
$ws.Range("E2").Value = '  +0.61%  '
$ws.Range("D3").Value = '2.265.04'
$ws.Range("E3").Value = '  -0.28%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("E5").Value = '  +4.29%  '
$ws.Range("D6").Value = '233.73'
$ws.Range("E6").Value = '  +0.92%  '
$ws.Range("D7").Value = '63.81'
$ws.Range("E7").Value = '  +0.43%  '
$ws.Range("D9").Value = '0.451'
$ws.Range("E9").Value = '  +4.42%  '
$ws.Range("D10").Value = '0.0977'
$ws.Range("E10").Value = '  -6.54%  '
$ws.Range("D11").Value = '57.89'
$ws.Range("E11").Value = '  +0.96%  '
$c = $ws.Range("D12")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '26.60'
$c.Style = $origStyle
$ws.Range("E12").Value = '  +2.69%  '
$ws.Range("E13").Value = '  +1.65%  '
$ws.Range("D14").Value = '2.600.49'
$ws.Range("E14").Value = '  -0.34%  '
$ws.Range("D15").Value = '15.64'
$ws.Range("E15").Value = '  -0.40%  '
$ws.Range("D16").Value = '6.15'
$ws.Range("E16").Value = '  +3.67%  '
$ws.Range("D17").Value = '0.842'
$ws.Range("E17").Value = '  +1.81%  '
$ws.Range("D18").Value = '2.269.09'
$ws.Range("E18").Value = '  +0.05%  '
$ws.Range("D19").Value = '43.841.29'
$ws.Range("E19").Value = '  +0.58%  '
$ws.Range("E20").Value = '  -2.12%  '
$ws.Range("E21").Value = '  +0.20%  '
$ws.Range("E22").Value = '  +0.68%  '
$ws.Range("D23").Value = '249.91'
$ws.Range("E23").Value = '  -0.07%  '
$ws.Range("E24").Value = '  -0.10%  '
$ws.Range("B25").Value = 'WEMIXToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D25").Value = '3.65'
$ws.Range("E25").Value = '  +30.52%  '
$ws.Range("B26").Value = 'PancakeSwap'
$ws.Range("C26").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D26").Value = '2.46'
$ws.Range("E26").Value = '  -2.21%  '
$ws.Range("E27").Value = '  -3.56%  '
$ws.Range("D28").Value = '9.89'
$ws.Range("E28").Value = '  -0.19%  '
$ws.Range("D29").Value = '174.06'
$ws.Range("E29").Value = '  +0.89%  '
$ws.Range("D30").Value = '21.98'
$ws.Range("E30").Value = '  +4.75%  '
$ws.Range("D31").Value = '0.137'
$ws.Range("E31").Value = '  -0.51%  '
$ws.Range("E32").Value = '  -0.61%  '
$ws.Range("E33").Value = '  +3.74%  '
$ws.Range("D34").Value = '5.01'
$ws.Range("E34").Value = '  +5.61%  '
$ws.Range("E35").Value = '  -0.23%  '
$ws.Range("E36").Value = '  -2.11%  '
$ws.Range("D37").Value = '3.71'
$ws.Range("E37").Value = '  -2.58%  '
$ws.Range("E38").Value = '  -5.69%  '
$ws.Range("E39").Value = '  -1.59%  '
$ws.Range("E40").Value = '  +2.86%  '
$ws.Range("E41").Value = '  +0.02%  '
$ws.Range("D42").Value = '8.76'
$ws.Range("E42").Value = '  +4.18%  '
$ws.Range("E43").Value = '  +3.41%  '
$ws.Range("D44").Value = '17.36'
$ws.Range("E44").Value = '  +0.21%  '
$ws.Range("D45").Value = '98.75'
$ws.Range("E45").Value = '  +1.10%  '
$ws.Range("B46").Value = 'Cronos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D46").Value = '0.0953'
$ws.Range("E46").Value = '  -0.74%  '
$ws.Range("B47").Value = 'TrustWalletToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D47").Value = '1.19'
$ws.Range("E47").Value = '  -1.11%  '
$ws.Range("E48").Value = '  +4.42%  '
$ws.Range("D49").Value = '1.458.37'
$ws.Range("E49").Value = '  -1.30%  '
$ws.Range("B50").Value = 'Celestia'
$ws.Range("C50").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D50").Value = '9.98'
$ws.Range("E50").Value = '  -3.95%  '
$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D51").Value = '2.31'
$ws.Range("E51").Value = '  -1.41%  '
